# Correcion a Diebold Mariano y revision de Cap1
#
# The DM-test summary table mis-assigned a few models/ratios; this fixes the
# "Modelo" labels, their significant-comparison ratios, and the associated
# Proporcion_Sig / Mejor_N_Calib / ECRPS_Mejor figures for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Column A (Modelo): fix model labels on rows 5, 6 and 8 ---
$ws.Cells.Item(5, 1).Value = "AREPD"
$ws.Cells.Item(6, 1).Value = "Block Bootstrapping"
$ws.Cells.Item(8, 1).Value = "DeepAR"

# --- Column B (Comparaciones_Significativas): update ratios on rows 2, 3 and 5 ---
$ws.Cells.Item(2, 2).Value = "3/10"
$ws.Cells.Item(3, 2).Value = "2/10"
$ws.Cells.Item(5, 2).Value = "0/10"

# --- Column C (Proporcion_Sig) ---
$ws.Cells.Item(2, 3).Value = 76.8
$ws.Cells.Item(3, 3).Value = 51.2
$ws.Cells.Item(5, 3).Value = 0

# --- Column D (Mejor_N_Calib) ---
$ws.Cells.Item(5, 4).Value = 40
$ws.Cells.Item(8, 4).Value = 200

# --- Column E (ECRPS_Mejor) ---
$ws.Cells.Item(5, 5).Value = 0.6733631690522695
$ws.Cells.Item(6, 5).Value = 0.6112845880987049
$ws.Cells.Item(8, 5).Value = 0.5901067674793075
